# A new, more recent price record was added for Albahaca at
# "Terminal La Palmera de La Serena". In the source sheet this new record
# sits at the top of the date-ordered data block (row 16, just under the
# already-present most-recent-ish rows), pushing every following data row
# down by one and growing the used range from A1:R120 to A1:R121.
#
# Net effect: insert one row at row 16 and populate it with the new
# observation; everything from the former row 16 onward shifts down
# automatically (row 120 -> row 121) with Excel handling the row-shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16; rows 16..120 shift to 17..121.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44819
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 100112052
$ws.Range("G16").Value = "Albahaca"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 3800
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = 3900
$ws.Range("N16").Value = "$/paquete"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 3900
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
